$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibition sheet): zero out "想去人数" (F) for all data rows ---
$wsExpo = $wb.Worksheets.Item("展览")
for ($r = 2; $r -le 14; $r++) {
    $wsExpo.Cells.Item($r, 6).Value = 0
}

# --- Sheet "演出" (performance sheet): zero out "想去人数" (F) for all data rows ---
$wsShow = $wb.Worksheets.Item("演出")
for ($r = 2; $r -le 5; $r++) {
    $wsShow.Cells.Item($r, 6).Value = 0
}

# --- Sheet "全部类型" (all-types sheet) ---
$wsAll = $wb.Worksheets.Item("全部类型")

# zero out "想去人数" (F) for the existing data rows (2..17)
for ($r = 2; $r -le 17; $r++) {
    $wsAll.Cells.Item($r, 6).Value = 0
}

# Duplicate row 17 (南宁·第二届北极光动漫展) into a freshly inserted row 18,
# pushing the old row 18 (南宁·万圣漫控嘉年华10) down to row 19.
$wsAll.Rows("17:17").Copy()
$wsAll.Rows("18:18").Insert()

# Fix up the sequence number in column A for the newly inserted row and the
# row that got shifted down.
$wsAll.Cells.Item(18, 1).Value = 17
$wsAll.Cells.Item(19, 1).Value = 18

# The shifted-down row (now row 19) keeps its own data but its "想去人数"
# must be reset to 0 as well.
$wsAll.Cells.Item(19, 6).Value = 0
